# Insert a new BOM row for a duplicated "10uF" capacitor (designator "    C41")
# right above the existing block that currently starts at row 71, pushing all
# following rows down by one (matches the xml diff: new row 71, everything
# else shifted +1, dimension A1:G148 -> A1:G149).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Insert a new blank row at position 71 (shifts old row71..148 down to 72..149)
$ws.Rows.Item(71).Insert()

# The new row 71 mirrors the other "10uF" capacitor rows (e.g. row 70, "    C30"),
# only the designator changes.
$ws.Range("A71").Value = "    C41"
$ws.Range("B71").Value = $ws.Range("B70").Value()
$ws.Range("C71").Value = $ws.Range("C70").Value()
$ws.Range("D71").Value = $ws.Range("D70").Value()
$ws.Range("E71").Value = $ws.Range("E70").Value()
$ws.Range("F71").Value = $ws.Range("F70").Value()
$ws.Range("G71").Value = 1

# Update the view: clear the scrolled-down top-left cell and move the
# selection/active cell to K10 (matches sheetView change in the diff).
$ws.Activate() | Out-Null
$ws.Range("K10").Select() | Out-Null
